# Daily attendance processing - 2026-01-16 01:41:30
# Normalizes the "Recorded By" column (G) so that the first two
# comma-separated entries are swapped (e.g. "user, System" -> "System, user").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $text = $cell.Text

    if ($text -and $text.Contains(",")) {
        $parts = $text -split ", "
        $firstIsSystem = $parts[0].Equals("System")
        if ($parts.Count -ge 2 -and -not $firstIsSystem) {
            $tmp = $parts[0]
            $parts[0] = $parts[1]
            $parts[1] = $tmp
            $cell.Value = $parts -join ", "
        }
    }
}
